$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 36; this shifts existing rows 36..100 down to 37..101
$ws.Rows.Item(36).Insert()

# Populate the newly inserted row 36 with the new weekly price record
$ws.Cells.Item(36, 1).Value = 10
$ws.Cells.Item(36, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(36, 3).Value = "La Araucanía"
$ws.Cells.Item(36, 4).Value = 45174
$ws.Cells.Item(36, 5).Value = 9
$ws.Cells.Item(36, 6).Value = 100112042
$ws.Cells.Item(36, 7).Value = "Locoto"
$ws.Cells.Item(36, 8).Value = "Sin especificar"
$ws.Cells.Item(36, 9).Value = "Primera"
$ws.Cells.Item(36, 10).Value = 100
$ws.Cells.Item(36, 11).Value = 2200
$ws.Cells.Item(36, 12).Value = 2200
$ws.Cells.Item(36, 13).Value = 2200
$ws.Cells.Item(36, 14).Value = "`$/kilo"
$ws.Cells.Item(36, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(36, 16).Value = 2200
$ws.Cells.Item(36, 17).Value = 1
$ws.Cells.Item(36, 18).Value = "Hortaliza"
